$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.973936579529778
$ws.Range("D2").Value = 7.729717733732095
$ws.Range("E2").Value = 12.72596180790718
$ws.Range("F2").Value = 44.58921583010547
$ws.Range("G2").Value = 57.8414135015265
$ws.Range("H2").Value = 19.80072001141792
$ws.Range("I2").Value = 17.51723746316345
$ws.Range("J2").Value = 10.109041464798
$ws.Range("N2").Value = 17.13588176063892

$ws.Range("B3").Value = 7.72610545600267
$ws.Range("D3").Value = 7.731607383498005
$ws.Range("E3").Value = 12.7347732422308
$ws.Range("F3").Value = 43.77798400427347
$ws.Range("G3").Value = 56.1295405628279
$ws.Range("H3").Value = 19.59404963810021
$ws.Range("I3").Value = 17.43545559361288
$ws.Range("J3").Value = 10.11490316061085
$ws.Range("N3").Value = 16.94880368210765

$ws.Range("B4").Value = 7.569968192699875
$ws.Range("D4").Value = 7.732971663349161
$ws.Range("E4").Value = 12.74289888234467
$ws.Range("F4").Value = 43.2888058150247
$ws.Range("G4").Value = 55.074807882116
$ws.Range("H4").Value = 19.47262906352406
$ws.Range("I4").Value = 17.38881070702369
$ws.Range("J4").Value = 10.12069478094539
$ws.Range("N4").Value = 16.83481638824602

$ws.Range("B5").Value = 7.505429865113573
$ws.Range("D5").Value = 7.733578766040622
$ws.Range("E5").Value = 12.74689087757529
$ws.Range("F5").Value = 43.09197246160743
$ws.Range("G5").Value = 54.64481053314584
$ws.Range("H5").Value = 19.42457083818628
$ws.Range("I5").Value = 17.37072111759998
$ws.Range("J5").Value = 10.12360443381637
$ws.Range("N5").Value = 16.78863424108461

$ws.Range("B6").Value = 7.494660968687292
$ws.Range("D6").Value = 7.733682659711756
$ws.Range("E6").Value = 12.74759479490273
$ws.Range("F6").Value = 43.05944780799597
$ws.Range("G6").Value = 54.57342058210207
$ws.Range("H6").Value = 19.4166779099519
$ws.Range("I6").Value = 17.36777342702472
$ws.Range("J6").Value = 10.12412071270919
$ws.Range("N6").Value = 16.78098333545053

$ws.Range("B7").Value = 7.569101375573497
$ws.Range("D7").Value = 7.732979644053993
$ws.Range("E7").Value = 12.74294996641711
$ws.Range("F7").Value = 43.2861407353301
$ws.Range("G7").Value = 55.06900854887076
$ws.Range("H7").Value = 19.47197512071688
$ws.Range("I7").Value = 17.38856299824766
$ws.Range("J7").Value = 10.12073179906159
$ws.Range("N7").Value = 16.83419240957985

$ws.Range("B8").Value = 7.889357710724791
$ws.Range("D8").Value = 7.730326842399959
$ws.Range("E8").Value = 12.72843522590959
$ws.Range("F8").Value = 44.30781285096198
$ws.Range("G8").Value = 57.25234060176211
$ws.Range("H8").Value = 19.72835136274256
$ws.Range("I8").Value = 17.48831115075864
$ws.Range("J8").Value = 10.11060648771349
$ws.Range("N8").Value = 17.07122325684961

$ws.Range("B9").Value = 8.482454446342578
$ws.Range("D9").Value = 7.726751176464113
$ws.Range("E9").Value = 12.72160251247386
$ws.Range("F9").Value = 46.37024980785385
$ws.Range("G9").Value = 61.47491567459172
$ws.Range("H9").Value = 20.27243875824717
$ws.Range("I9").Value = 17.71126666152004
$ws.Range("J9").Value = 10.10822160243079
$ws.Range("N9").Value = 17.54092831277455

$ws.Range("B10").Value = 8.892805469390057
$ws.Range("D10").Value = 7.725127013220942
$ws.Range("E10").Value = 12.72987491738833
$ws.Range("F10").Value = 47.90568748517494
$ws.Range("G10").Value = 64.50414412291448
$ws.Range("H10").Value = 20.69436942026452
$ws.Range("I10").Value = 17.89037526883945
$ws.Range("J10").Value = 10.11721207091131
$ws.Range("N10").Value = 17.88621396499382

$ws.Range("B11").Value = 9.073239364152281
$ws.Range("D11").Value = 7.724608544374606
$ws.Range("E11").Value = 12.73654248767226
$ws.Range("F11").Value = 48.6053310187205
$ws.Range("G11").Value = 65.85944379355729
$ws.Range("H11").Value = 20.89045543231269
$ws.Range("I11").Value = 17.97487022382016
$ws.Range("J11").Value = 10.12365035974141
$ws.Range("N11").Value = 18.04275315944065

$ws.Range("B12").Value = 9.140615808878721
$ws.Range("D12").Value = 7.724444128018039
$ws.Range("E12").Value = 12.73948595990004
$ws.Range("F12").Value = 48.87017583043146
$ws.Range("G12").Value = 66.368886747877
$ws.Range("H12").Value = 20.96524579847746
$ws.Range("I12").Value = 18.00727319227334
$ws.Range("J12").Value = 10.12642695060378
$ws.Range("N12").Value = 18.10190618682536

$ws.Range("B13").Value = 9.126148118074379
$ws.Range("D13").Value = 7.724478114804402
$ws.Range("E13").Value = 12.73883340025739
$ws.Range("F13").Value = 48.81314504986548
$ws.Range("G13").Value = 66.25934484517246
$ws.Range("H13").Value = 20.9491154314576
$ws.Range("I13").Value = 18.00027699228538
$ws.Range("J13").Value = 10.12581389385634
$ws.Range("N13").Value = 18.08917283935262

$ws.Range("B14").Value = 9.078801786399715
$ws.Range("D14").Value = 7.724594376744808
$ws.Range("E14").Value = 12.73677625479573
$ws.Range("F14").Value = 48.62712332969024
$ws.Range("G14").Value = 65.90143448747506
$ws.Range("H14").Value = 20.89659806435153
$ws.Range("I14").Value = 17.97752804593555
$ws.Range("J14").Value = 10.12387200243344
$ws.Range("N14").Value = 18.04762249193309

$ws.Range("B15").Value = 9.049675604494462
$ws.Range("D15").Value = 7.724669753947895
$ws.Range("E15").Value = 12.73557073299515
$ws.Range("F15").Value = 48.51315967612338
$ws.Range("G15").Value = 65.68169744160305
$ws.Range("H15").Value = 20.86449779800469
$ws.Range("I15").Value = 17.96364578513091
$ws.Range("J15").Value = 10.12272664843227
$ws.Range("N15").Value = 18.02215401479317

$ws.Range("B16").Value = 8.880883512103265
$ws.Range("D16").Value = 7.72516535002237
$ws.Range("E16").Value = 12.72949770781503
$ws.Range("F16").Value = 47.85996489177818
$ws.Range("G16").Value = 64.41507003753819
$ws.Range("H16").Value = 20.68163311876464
$ws.Range("I16").Value = 17.88491171931946
$ws.Range("J16").Value = 10.11683864795549
$ws.Range("N16").Value = 17.87596887941635

$ws.Range("B17").Value = 8.775697920554402
$ws.Range("D17").Value = 7.725525983146547
$ws.Range("E17").Value = 12.72651675488256
$ws.Range("F17").Value = 47.45935932342474
$ws.Range("G17").Value = 63.63183684729849
$ws.Range("H17").Value = 20.5704716703448
$ws.Range("I17").Value = 17.83736449730079
$ws.Range("J17").Value = 10.11382872043847
$ws.Range("N17").Value = 17.78611779882747

$ws.Range("B18").Value = 8.714613687309594
$ws.Range("D18").Value = 7.725754137937617
$ws.Range("E18").Value = 12.72507552557142
$ws.Range("F18").Value = 47.22906635946156
$ws.Range("G18").Value = 63.17923164279156
$ws.Range("H18").Value = 20.50692855796396
$ws.Range("I18").Value = 17.81030229408064
$ws.Range("J18").Value = 10.11231850326613
$ws.Range("N18").Value = 17.7343904948633

$ws.Range("B19").Value = 8.693832984641853
$ws.Range("D19").Value = 7.725834940556858
$ws.Range("E19").Value = 12.72463446245567
$ws.Range("F19").Value = 47.15112272176361
$ws.Range("G19").Value = 63.02564133938778
$ws.Range("H19").Value = 20.48548342974222
$ws.Range("I19").Value = 17.80118939117584
$ws.Range("J19").Value = 10.11184509972871
$ws.Range("N19").Value = 17.71686987832946

$ws.Range("B20").Value = 8.786955979029527
$ws.Range("D20").Value = 7.725485446234032
$ws.Range("E20").Value = 12.72680578544053
$ws.Range("F20").Value = 47.50199349396804
$ws.Range("G20").Value = 63.71543563327825
$ws.Range("H20").Value = 20.58226464743871
$ws.Range("I20").Value = 17.84239660774302
$ws.Range("J20").Value = 10.11412625151573
$ws.Range("N20").Value = 17.79568785271756

$ws.Range("B21").Value = 9.0927347269932
$ws.Range("D21").Value = 7.724559359728069
$ws.Range("E21").Value = 12.7373691210164
$ws.Range("F21").Value = 48.6817670345036
$ws.Range("G21").Value = 66.00666786147329
$ws.Range("H21").Value = 20.91200959515086
$ws.Range("I21").Value = 17.9841991494933
$ws.Range("J21").Value = 10.12443318921646
$ws.Range("N21").Value = 18.05983060728496

$ws.Range("B22").Value = 9.287021218569393
$ws.Range("D22").Value = 7.724140255687939
$ws.Range("E22").Value = 12.74671310399783
$ws.Range("F22").Value = 49.45214975722511
$ws.Range("G22").Value = 67.48188650198634
$ws.Range("H22").Value = 21.13061830241064
$ws.Range("I22").Value = 18.07923203987225
$ws.Range("J22").Value = 10.13314293819171
$ws.Range("N22").Value = 18.23171317604656

$ws.Range("B23").Value = 9.183851410766957
$ws.Range("D23").Value = 7.724346829107239
$ws.Range("E23").Value = 12.74150251249817
$ws.Range("F23").Value = 49.0411268991405
$ws.Range("G23").Value = 66.69672593583009
$ws.Range("H23").Value = 21.01367872353836
$ws.Range("I23").Value = 18.02830472145925
$ws.Range("J23").Value = 10.12831357741493
$ws.Range("N23").Value = 18.1400605029821

$ws.Range("B24").Value = 8.78186811180856
$ws.Range("D24").Value = 7.725503708133549
$ws.Range("E24").Value = 12.72667426583628
$ws.Range("F24").Value = 47.48271850152633
$ws.Range("G24").Value = 63.67764779024028
$ws.Range("H24").Value = 20.57693190082308
$ws.Range("I24").Value = 17.84012073718849
$ws.Range("J24").Value = 10.11399105181379
$ws.Range("N24").Value = 17.79136144795663

$ws.Range("B25").Value = 8.326196084739975
$ws.Range("D25").Value = 7.727543224494022
$ws.Range("E25").Value = 12.72112336097148
$ws.Range("F25").Value = 45.80763838381693
$ws.Range("G25").Value = 60.34292603907729
$ws.Range("H25").Value = 20.12112573215001
$ws.Range("I25").Value = 17.64816122966608
$ws.Range("J25").Value = 10.10698609646884
$ws.Range("N25").Value = 17.41362653806352
